{"js": "// Update the date title and the 25 division-problem answers in the table.\n// The table has 20 rows total but only every 4th row (0, 4, 8, 12, 16)\n// actually contains the 5 columns of division problems; the rows in\n// between are blank spacer rows.\n\nconst title = context.document.body.paragraphs.getFirst();\ntitle.load(\"text\");\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst newTitle = \"2025-06-19 Thursday\";\n\n// Replace the title text (first paragraph) while preserving its run\n// formatting: use search on the paragraph range so formatting of the\n// existing run is kept.\nconst titleResults = title.search(\"2025-06-18 Wednesday\", { matchCase: true });\ntitleResults.load(\"text\");\nawait context.sync();\nif (titleResults.items.length > 0) {\n  titleResults.items[0].insertText(newTitle, Word.InsertLocation.replace);\n} else {\n  // Fallback: if the exact text wasn't found (already changed, etc.)\n  // just overwrite the whole paragraph text.\n  title.insertText(newTitle, Word.InsertLocation.replace);\n}\n\n// Grid of new values, row-major, 5 columns per row. These correspond to\n// the table rows that actually hold data: 0, 4, 8, 12, 16.\nconst newValues = [\n  [\"55\u00f72=27, 1\", \"21\u00f72=10, 1\", \"89\u00f74=22, 1\", \"80\u00f73=26, 2\", \"88\u00f79=9, 7\"],\n  [\"41\u00f73=13, 2\", \"27\u00f74=6, 3\", \"93\u00f75=18, 3\", \"36\u00f75=7, 1\", \"30\u00f77=4, 2\"],\n  [\"30\u00f76=5, 0\", \"66\u00f77=9, 3\", \"69\u00f72=34, 1\", \"94\u00f79=10, 4\", \"27\u00f76=4, 3\"],\n  [\"89\u00f74=22, 1\", \"59\u00f77=8, 3\", \"66\u00f78=8, 2\", \"81\u00f79=9, 0\", \"50\u00f73=16, 2\"],\n  [\"88\u00f77=12, 4\", \"19\u00f73=6, 1\", \"95\u00f77=13, 4\", \"23\u00f73=7, 2\", \"85\u00f72=42, 1\"],\n];\n\nconst dataRowIndices = [0, 4, 8, 12, 16];\n\nfor (let r = 0; r < dataRowIndices.length; r++) {\n  const tableRowIndex = dataRowIndices[r];\n  for (let c = 0; c < 5; c++) {\n    const cell = table.getCell(tableRowIndex, c);\n    const cellBody = cell.body;\n    cellBody.load(\"text\");\n    // eslint-disable-next-line no-await-in-loop\n    await context.sync();\n\n    const oldText = cellBody.text.trim();\n    const newText = newValues[r][c];\n\n    const found = cellBody.search(oldText, { matchCase: true });\n    found.load(\"text\");\n    // eslint-disable-next-line no-await-in-loop\n    await context.sync();\n\n    if (found.items.length > 0) {\n      found.items[0].insertText(newText, Word.InsertLocation.replace);\n    } else {\n      cellBody.insertText(newText, Word.InsertLocation.replace);\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date title and the 25 division-problem answers in the table.\n# Every text run that carries content in the document changes, and every\n# \"old\" value is unique within the document, so a plain exact-text\n# Find/Replace over the whole document body is safe and unambiguous.\n\n$d = $word.ActiveDocument\n\nfunction Replace-ExactText($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\n# Title paragraph (the worksheet date).\nReplace-ExactText '2025-06-18 Wednesday' '2025-06-19 Thursday'\n\n# Table of division problems (5 rows x 5 columns of \"a\u00f7b=c, d\" answers).\nReplace-ExactText '28\u00f72=14, 0' '55\u00f72=27, 1'\nReplace-ExactText '36\u00f79=4, 0' '21\u00f72=10, 1'\nReplace-ExactText '61\u00f75=12, 1' '89\u00f74=22, 1'\nReplace-ExactText '41\u00f79=4, 5' '80\u00f73=26, 2'\nReplace-ExactText '61\u00f79=6, 7' '88\u00f79=9, 7'\n\nReplace-ExactText '56\u00f75=11, 1' '41\u00f73=13, 2'\nReplace-ExactText '79\u00f78=9, 7' '27\u00f74=6, 3'\nReplace-ExactText '79\u00f74=19, 3' '93\u00f75=18, 3'\nReplace-ExactText '14\u00f75=2, 4' '36\u00f75=7, 1'\nReplace-ExactText '42\u00f77=6, 0' '30\u00f77=4, 2'\n\nReplace-ExactText '72\u00f78=9, 0' '30\u00f76=5, 0'\nReplace-ExactText '22\u00f78=2, 6' '66\u00f77=9, 3'\nReplace-ExactText '16\u00f77=2, 2' '69\u00f72=34, 1'\nReplace-ExactText '76\u00f77=10, 6' '94\u00f79=10, 4'\nReplace-ExactText '10\u00f73=3, 1' '27\u00f76=4, 3'\n\nReplace-ExactText '86\u00f78=10, 6' '89\u00f74=22, 1'\nReplace-ExactText '31\u00f73=10, 1' '59\u00f77=8, 3'\nReplace-ExactText '61\u00f76=10, 1' '66\u00f78=8, 2'\nReplace-ExactText '58\u00f73=19, 1' '81\u00f79=9, 0'\nReplace-ExactText '25\u00f73=8, 1' '50\u00f73=16, 2'\n\nReplace-ExactText '58\u00f76=9, 4' '88\u00f77=12, 4'\nReplace-ExactText '39\u00f75=7, 4' '19\u00f73=6, 1'\nReplace-ExactText '34\u00f79=3, 7' '95\u00f77=13, 4'\nReplace-ExactText '60\u00f74=15, 0' '23\u00f73=7, 2'\nReplace-ExactText '96\u00f77=13, 5' '85\u00f72=42, 1'\n"}
